# GasMileageTestData.xlsx — record the 13.33 mpg result for row 2
# (Expected column E and Actual column F both read "13.33").
#
# Apache POI's Cell#setCellValue(String) stores a numeric-looking value as a
# plain text/shared-string cell (no special number formatting). To reproduce
# that via Excel COM automation — where typing "13.33" into a General cell
# is normally auto-coerced to a number — we briefly mark the cells as Text
# before assigning the value, then clear the formatting again so the cells
# are left without any explicit style, just like the source edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targets = $ws.Range("E2:F2")
$targets.NumberFormat = "@"
$ws.Range("E2").Value = "13.33"
$ws.Range("F2").Value = "13.33"
$targets.ClearFormats()
